# Updated cryptos list on Wed Aug 23 21:07:24 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row (rows 2-51) on the active sheet with freshly scraped
# values. Prices are kept as literal text (matching the sheet's existing
# "26.635.28"-style formatted strings) rather than numbers, so rows whose
# new price would otherwise be auto-parsed by Excel as a plain number
# (DIsText = $true) have their cell's NumberFormat forced to "@" (Text)
# right before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "26.630.04"; DIsText = $false; E = "  +2.74%  " },
    @{ Row = 3; D = "1.686.58"; DIsText = $false; E = "  +3.31%  " },
    @{ Row = 4; D = $null; DIsText = $false; E = "  -0.04%  " },
    @{ Row = 5; D = "217.14"; DIsText = $true; E = "  +3.61%  " },
    @{ Row = 6; D = "0.5333"; DIsText = $true; E = "  +2.62%  " },
    @{ Row = 7; D = $null; DIsText = $false; E = "  -0.05%  " },
    @{ Row = 8; D = "0.2680"; DIsText = $true; E = "  +4.49%  " },
    @{ Row = 9; D = "0.06424"; DIsText = $true; E = "  +3.10%  " },
    @{ Row = 10; D = "21.66"; DIsText = $true; E = "  +6.91%  " },
    @{ Row = 11; D = "0.07799"; DIsText = $true; E = "  +3.03%  " },
    @{ Row = 12; D = "1.699.17"; DIsText = $false; E = "  +4.53%  " },
    @{ Row = 13; D = "4.497"; DIsText = $true; E = "  +3.38%  " },
    @{ Row = 14; D = "0.5601"; DIsText = $true; E = "  +3.41%  " },
    @{ Row = 15; D = "0.0₅8416"; DIsText = $false; E = "  +5.99%  " },
    @{ Row = 16; D = "66.04"; DIsText = $true; E = "  +2.34%  " },
    @{ Row = 17; D = "26.672.86"; DIsText = $false; E = "  +2.89%  " },
    @{ Row = 18; D = $null; DIsText = $false; E = "  -0.09%  " },
    @{ Row = 19; D = "4.797"; DIsText = $true; E = "  +3.88%  " },
    @{ Row = 20; D = "195.26"; DIsText = $true; E = "  +5.84%  " },
    @{ Row = 21; D = "10.41"; DIsText = $true; E = "  +3.97%  " },
    @{ Row = 22; D = "6.379"; DIsText = $true; E = "  +5.10%  " },
    @{ Row = 23; D = "1.003"; DIsText = $true; E = "  -0.05%  " },
    @{ Row = 24; D = "144.13"; DIsText = $true; E = "  -1.10%  " },
    @{ Row = 25; D = "0.1282"; DIsText = $true; E = "  +6.28%  " },
    @{ Row = 26; D = "7.470"; DIsText = $true; E = "  +1.73%  " },
    @{ Row = 27; D = "16.27"; DIsText = $true; E = "  +5.06%  " },
    @{ Row = 28; D = "1.433"; DIsText = $true; E = "  +4.69%  " },
    @{ Row = 29; D = "0.06148"; DIsText = $true; E = "  +3.70%  " },
    @{ Row = 30; D = "1.279"; DIsText = $true; E = "  +3.02%  " },
    @{ Row = 31; D = "3.607"; DIsText = $true; E = "  +7.54%  " },
    @{ Row = 32; D = "3.466"; DIsText = $true; E = "  +3.50%  " },
    @{ Row = 33; D = "1.699"; DIsText = $true; E = "  +5.76%  " },
    @{ Row = 34; D = "1.013"; DIsText = $true; E = "  +4.38%  " },
    @{ Row = 35; D = "2.423"; DIsText = $true; E = "  +1.71%  " },
    @{ Row = 36; D = "2.793"; DIsText = $true; E = "  +2.07%  " },
    @{ Row = 37; D = "0.5741"; DIsText = $true; E = "  -0.55%  " },
    @{ Row = 38; D = "0.01647"; DIsText = $true; E = "  +3.22%  " },
    @{ Row = 39; D = "6.023"; DIsText = $true; E = "  +6.96%  " },
    @{ Row = 40; D = "1.069.74"; DIsText = $false; E = "  +5.36%  " },
    @{ Row = 41; D = "0.8639"; DIsText = $true; E = "  +3.18%  " },
    @{ Row = 42; D = "1.002"; DIsText = $true; E = $null },
    @{ Row = 43; D = "100.38"; DIsText = $true; E = "  +0.82%  " },
    @{ Row = 44; D = "1.837.87"; DIsText = $false; E = "  +3.04%  " },
    @{ Row = 45; D = $null; DIsText = $false; E = "  +3.41%  " },
    @{ Row = 46; D = "57.25"; DIsText = $true; E = "  +5.53%  " },
    @{ Row = 47; D = "8.217"; DIsText = $true; E = "  +3.39%  " },
    @{ Row = 48; D = "1.003"; DIsText = $true; E = "  +0.37%  " },
    @{ Row = 49; D = "0.05228"; DIsText = $true; E = "  +1.07%  " },
    @{ Row = 50; D = "6.086"; DIsText = $true; E = "  +5.31%  " },
    @{ Row = 51; D = "0.4242"; DIsText = $true; E = "  +0.41%  " }
)

foreach ($row in $updates) {
    if ($null -ne $row.D) {
        $cell = $ws.Cells.Item($row.Row, 4)
        if ($row.DIsText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row.D
    }
    if ($null -ne $row.E) {
        $ws.Cells.Item($row.Row, 5).Value = $row.E
    }
}
